# A new weekly price observation is inserted as row 81 in the Frambuesa
# (raspberry) price series. All existing rows from 81 downward shift down
# by one row (to 82:201), and the new row 81 is populated with the new
# observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 81 - shifts rows 81:200 down to 82:201.
$ws.Rows("81").Insert()

# Populate the new row 81 with the new observation. The columns that are
# constant across the whole sheet (Mercado ID, Mercado, Región, Codreg,
# Tipo, Producto ID, Producto, Categoría ID, Categoría, Variedad) are
# copied from the row immediately below (now row 82, the original row 81).
$ws.Range("A81").Value = 6
$ws.Range("B81").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C81").Value = "Metropolitana"
$ws.Range("D81").Value = 44671
$ws.Range("E81").Value = 13
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100101
$ws.Range("H81").Value = "Berries"
$ws.Range("I81").Value = 100101004
$ws.Range("J81").Value = "Frambuesa"
$ws.Range("K81").Value = "Sin especificar"
$ws.Range("L81").Value = "Primera"
$ws.Range("M81").Value = 200
$ws.Range("N81").Value = 8000
$ws.Range("O81").Value = 8000
$ws.Range("P81").Value = 8000
$ws.Range("Q81").Value = "$/bandeja 2 kilos"
$ws.Range("R81").Value = "Región del Maule"
$ws.Range("S81").Value = 4000
$ws.Range("T81").Value = 2
